$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2:O73").Value = "2022-08-17 20:59:50"
